# Generate Report for Handoff
# Update the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" values
# for the f75bf15b-8994-4924-9b93-459ee87e195c entry on all three sheets, as
# part of a fresh handoff report generation.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G7").Value = "2016-08-20 12:45:30"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H7").Value = "2016-08-20 12:45:25"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H7").Value = "2016-08-20 12:45:30"
